# Trade #74 closed at 2026-02-17 08:57:58 - unknown UNKNOWN +0.000%
#
# This script updates the live trading results workbook to record the
# closing of trade #74:
#   - Summary sheet: bump Total Trades (B6) and refresh Win Rate % (B9)
#   - Strategy Status sheet: bump MarketMaking Trades (D4) and refresh its
#     Win Rate % (G4)
#   - All Trades sheet: append the new trade as row 75
#   - MarketMaking sheet: append the same trade as row 75

$wb = $excel.ActiveWorkbook

# ---- Summary sheet ----
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 74
$summary.Range("B9").Value = 41.89

# ---- Strategy Status sheet ----
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 74
$status.Range("G4").Value = 41.89

# ---- New trade row data (shared by "All Trades" and "MarketMaking") ----
$tradeNum = 74
$tradeDate = "2026-02-17"
$tradeTime = "08:57:52"
$strategy = "MarketMaking"
$side = "DOWN"
$entryPrice = 0.89
$exitPrice = 0.89
$tradeStatus = "CLOSED"
$plPct = 0
$plDollar = 0
$capitalAfter = 100.57
$entrySlippage = 0
$exitSlippage = 0
$confidence = 0.6
$entryReason = "Normal spread capture: 19600 bps"
$exitReason = "early_exit"
$duration = 0.13

$row = 75

function Set-TradeRow($ws, $row) {
    $ws.Cells.Item($row, 1).Value = $tradeNum

    # Column B holds a date-like string ("2026-02-17"). Force it to be
    # stored as text (matching the rest of the sheet) instead of being
    # auto-converted into a date serial number, then drop the temporary
    # number format so the cell is left with no special styling, just
    # like every other cell in the sheet.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = $tradeDate
    $ws.Cells.Item($row, 2).ClearFormats()

    $ws.Cells.Item($row, 3).Value = $tradeTime
    $ws.Cells.Item($row, 4).Value = $strategy
    $ws.Cells.Item($row, 5).Value = $side
    $ws.Cells.Item($row, 6).Value = $entryPrice
    $ws.Cells.Item($row, 7).Value = $exitPrice
    $ws.Cells.Item($row, 8).Value = $tradeStatus
    $ws.Cells.Item($row, 9).Value = $plPct
    $ws.Cells.Item($row, 10).Value = $plDollar
    $ws.Cells.Item($row, 11).Value = $capitalAfter
    $ws.Cells.Item($row, 12).Value = $entrySlippage
    $ws.Cells.Item($row, 13).Value = $exitSlippage
    $ws.Cells.Item($row, 14).Value = $confidence
    $ws.Cells.Item($row, 15).Value = $entryReason
    $ws.Cells.Item($row, 16).Value = $exitReason
    $ws.Cells.Item($row, 17).Value = $duration
}

# ---- All Trades sheet ----
$allTrades = $wb.Worksheets.Item("All Trades")
Set-TradeRow $allTrades $row

# ---- MarketMaking sheet ----
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Set-TradeRow $marketMaking $row
